$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44261.54823399689
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44261.52689609954
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44261.50553879629
}
